$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2: "Lala" -> "Matheus"
$ws.Range("A2").Value = "Matheus"

# B2: numeric 321 -> text "123" (store as text, not a number)
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "123"
$ws.Range("B2").ClearFormats()

# C2: date serial 37289 (2002-02-02 with time) -> 37120 (2001-08-17, date only)
$ws.Range("C2").Value = (Get-Date -Year 2001 -Month 8 -Day 17 -Hour 0 -Minute 0 -Second 0)
$ws.Range("C2").NumberFormat = "YYYY-MM-DD"
